# The commit swaps the presentation's theme from the custom "Integral"
# palette to the stock "Office" palette (ppt/theme/theme1.xml, the theme
# wired to the slide master). Drive this the same way a user would from
# the Design tab / Variants gallery: push the 12 theme colour slots on
# the slide master's theme to the standard Office RGB values.
#
# ThemeColorScheme.Item(n) slot order is the usual DrawingML clrScheme
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. The .RGB setter
# takes a COLORREF (0xBBGGRR), so each target "RRGGBB" hex value below is
# byte-reversed before being assigned.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000  # dk1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444  # dk2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED  # accent2  -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$cs.Item(9).RGB  = 0xC47244  # accent5  -> 4472C4
$cs.Item(10).RGB = 0x47AD70  # accent6  -> 70AD47
$cs.Item(11).RGB = 0xC16305  # hlink    -> 0563C1
$cs.Item(12).RGB = 0x724F95  # folHlink -> 954F72
